# Atualização automática: 2025-09-01 09:00:26
# Applies the row-level data corrections to the detections sheet:
#   - Row 18: update first-detection image filename and refine the bounding
#     box / confidence for that detection.
#   - Rows 21/22: the two "mosca" detections for 2025-08-24 swap places
#     (Fly_ID + downstream bbox/confidence), and the record now sitting in
#     row 21 gets a corrected image/coords/confidence.
#
# Columns I (coords) and J (confidence) are plain text cells in the source
# workbook (t="inlineStr"), even though their contents look numeric. Setting
# .Value on a numeric-looking string makes Excel coerce it into a real
# number (and tags the cell with a "@" text format), so for those cells we
# briefly force Text format, assign the literal text, then clear the format
# back off so no stray style survives - matching the original "plain text,
# default style" shape of these cells.

function Set-TextValue($Cell, $Text) {
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18 ---------------------------------------------------------------
$ws.Range("D18").Value = "image_20250808221835_ppp0.jpg"
Set-TextValue $ws.Range("I18") "1182,405,1231,455"
Set-TextValue $ws.Range("J18") "0.76"

# --- Row 21 (now holds the former row-22 record, with refreshed detection
#     details) -------------------------------------------------------------
$ws.Range("A21").Value = "a2ea21b8-7dce-4e6a-be35-4edaddca5896"
$ws.Range("D21").Value = "image_20250824092407_ppp0.jpg"
Set-TextValue $ws.Range("I21") "1002,789,1039,825"
Set-TextValue $ws.Range("J21") "0.64"

# --- Row 22 (now holds the former row-21 record) ---------------------------
$ws.Range("A22").Value = "66efa766-1456-4beb-b92a-0615a2fc41bb"
Set-TextValue $ws.Range("I22") "1272,293,1315,331"
Set-TextValue $ws.Range("J22") "0.69"
